$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New rows of data to append: date (serial), epidemiological_week,
# last_available_confirmed, last_available_deaths, new_confirmed, new_deaths
$newRows = @(
    @(44746, 0, 331286, 6360, 83, 0),
    @(44747, 0, 331729, 6361, 443, 1),
    @(44748, 0, 332454, 6363, 725, 2),
    @(44749, 0, 332978, 6364, 524, 1),
    @(44750, 0, 333606, 6366, 628, 2)
)

$startRow = 91
$r = $startRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$ws.Range("F94").Select()
